$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Part 1: "Frontend environment: React -> Vite." becomes
#         "Frontend environment: React -> React Vite."
#         (the bold run is conceptually split into three runs with
#          identical run-properties; a plain find/replace reproduces
#          the same visible text & formatting)
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "React -> Vite.", $false, $false, $false, $false, $false,
    $true, 1, $false, "React -> React Vite.", 2) | Out-Null

# -----------------------------------------------------------------
# Part 2: add a brand-new dated journal entry after the
#         "...mint middleware." paragraph (and before the last,
#         empty paragraph of the document):
#
#   2024.08.01.                                   (bold)
#   Frontend routing átalakítása és finomítása (URL). További
#   hibajavítások, optimalizálások. Adatbázis optimalizáció
#   validációnak megfelelően.                     (regular)
# -----------------------------------------------------------------

# Grab a fully-formatted "plain" (non-bold) run already in the
# document so the new body text can reuse its exact rPr (font,
# size, language) without having to rebuild it property-by-property.
$sourcePara = $d.Paragraphs(98)
$plainTailLen = "mint middleware.".Length
$plainRefStart = $sourcePara.Range.End - $plainTailLen
$plainRef = $d.Range($plainRefStart, $sourcePara.Range.End)
$plainRefText = $plainRef.Text
$plainFormatted = $plainRef.FormattedText

# The very last paragraph in the document is empty; its paragraph
# mark already carries bold formatting (Times New Roman, 13pt,
# hu-HU) which is exactly what the new date heading needs, so we
# insert in front of it.
$finalPara = $d.Paragraphs($d.Paragraphs.Count)
$insertStart = $finalPara.Range.Start

$dateText = "2024.08.01."
$placeholder = "XJournalBodyPlaceholderX"
$newParaText = $dateText + [char]11 + $placeholder

$insertPoint = $d.Range($insertStart, $insertStart)
$insertPoint.InsertBefore($newParaText + [char]13)

# Re-point the placeholder body text to non-bold formatting, copied
# from the reference plain run above.
$placeholderStart = $insertStart + $dateText.Length + 1
$placeholderEnd = $placeholderStart + $placeholder.Length
$placeholderRange = $d.Range($placeholderStart, $placeholderEnd)
$placeholderRange.FormattedText = $plainFormatted

# Swap the placeholder text for the real paragraph body (using the
# Range.Text setter keeps the rPr that FormattedText just applied).
$bodyText = "Frontend routing átalakítása és finomítása (URL). " + `
    "További hibajavítások, optimalizálások. Adatbázis optimalizáció " + `
    "validációnak megfelelően."
$bodyRange = $d.Range($placeholderStart, $placeholderStart + $plainRefText.Length)
$bodyRange.Text = $bodyText
